{"js": "// The document opens with a leading paragraph that contains a single\n// inline picture (the \"Empty red chairs...\" image). The edit removes\n// the picture but keeps the (now empty) paragraph in place, right\n// before the customer-info table.\nconst body = context.document.body;\nconst pictures = body.inlinePictures;\npictures.load(\"items\");\nawait context.sync();\n\nfor (const picture of pictures.items) {\n  picture.delete();\n}\nawait context.sync();\n", "ps1": "# The document opens with a leading paragraph that contains a single\n# inline picture (the \"Empty red chairs...\" image). Remove the picture\n# but leave the (now empty) paragraph in place, right before the\n# customer-info table.\n$d = $word.ActiveDocument\n\nfor ($i = $d.InlineShapes.Count; $i -ge 1; $i--) {\n    $d.InlineShapes.Item($i).Delete()\n}\n"}
